# Generate Report for Handback
#
# The nightly localization-status report is regenerated: the two files
# ("179e0ac6...md" / "89e4f584...md") have now been handed back from
# localization, so:
#   - overall status flips from "Ready for handoff" to
#     "Handed back: in sync with en-US" everywhere it appears
#     (Overview!E2:F3 and the per-locale Status column C on zh-cn/de-de)
#   - the per-locale rows gain the "Latest Target File" (md, hyperlinked)
#     and "Latest Handback File" (xlf) values that were previously blank
#   - the "Latest Handback DateTime" column gets populated with the
#     actual handback timestamps

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$mdFile1 = "179e0ac6-954a-424c-80fa-b0d0b67df686.md"
$mdFile2 = "89e4f584-c1e2-4ab6-a63d-3f823a5d5245.md"

$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8d2ce6f8e54d895b65f66b331ca9df72a0bc23f5/e2e/$mdFile1"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8d2ce6f8e54d895b65f66b331ca9df72a0bc23f5/e2e/$mdFile2"

$zhXlf1 = "179e0ac6-954a-424c-80fa-b0d0b67df686.180c74fb69f314232aae9ef5d110a6e0d875ea0e.zh-cn.xlf"
$zhXlf2 = "89e4f584-c1e2-4ab6-a63d-3f823a5d5245.d752b22f459b25070335ed23e59b9942d5db3a03.zh-cn.xlf"
$deXlf1 = "179e0ac6-954a-424c-80fa-b0d0b67df686.180c74fb69f314232aae9ef5d110a6e0d875ea0e.de-de.xlf"
$deXlf2 = "89e4f584-c1e2-4ab6-a63d-3f823a5d5245.d752b22f459b25070335ed23e59b9942d5db3a03.de-de.xlf"

$zhHandbackTime = "2016-09-05 10:40:54"
$deHandbackTime = "2016-09-05 10:41:09"

# Column width that matches the widened columns (custom width "29.9777..."
# in the saved file); 30 characters is the closest Excel lets us land on
# through the ColumnWidth property (Excel quantizes column widths to
# whole pixels).
$wideWidth = 29.17
# Column width that matches the widened "Latest Target File" / "Latest
# Handback File" columns (40 characters, same as the other fixed-width
# 40-char columns already in the sheet).
$fullWidth = 39.17

# ---------------------------------------------------------------------
# Overview sheet: status text + widen the zh-cn / de-de status columns
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("E2").Value = $statusNew
$ov.Range("F2").Value = $statusNew
$ov.Range("E3").Value = $statusNew
$ov.Range("F3").Value = $statusNew

$ov.Columns.Item(5).ColumnWidth = $wideWidth
$ov.Columns.Item(6).ColumnWidth = $wideWidth

# ---------------------------------------------------------------------
# Per-locale sheets (zh-cn, de-de): status text, target/handback files,
# handback datetime, widened columns, and hyperlinks on the newly
# populated "Latest Target File" cells.
# ---------------------------------------------------------------------
function Update-LocaleSheet {
    param($sheetName, $xlf1, $xlf2, $handbackTime)
    $ws = $wb.Worksheets.Item($sheetName)

    # Status column (C) flips to the handed-back status too.
    $ws.Range("C2").Value = $statusNew
    $ws.Range("C3").Value = $statusNew

    # "Latest Target File" (I) - the localized markdown file, hyperlinked
    # the same way the "Source File Name" (A) column already is.
    $ws.Range("I2").Value = $mdFile1
    $ws.Range("I3").Value = $mdFile2

    # "Latest Handback File" (J) - the generated xliff for this locale.
    $ws.Range("J2").Value = $xlf1
    $ws.Range("J3").Value = $xlf2

    # "Latest Handback DateTime" (K).
    $ws.Range("K2").Value = $handbackTime
    $ws.Range("K3").Value = $handbackTime

    # Hyperlinks on the newly-filled target-file cells, mirroring A2/A3.
    $ws.Hyperlinks.Add($ws.Range("I2"), $mdUrl1, "", "", $mdFile1)
    $ws.Hyperlinks.Add($ws.Range("I3"), $mdUrl2, "", "", $mdFile2)

    # Widen Status (C) and the now-longer Target/Handback file columns
    # (I, J) so the new long filenames are readable.
    $ws.Columns.Item(3).ColumnWidth = $wideWidth
    $ws.Columns.Item(9).ColumnWidth = $fullWidth
    $ws.Columns.Item(10).ColumnWidth = $fullWidth
}

Update-LocaleSheet "zh-cn" $zhXlf1 $zhXlf2 $zhHandbackTime
Update-LocaleSheet "de-de" $deXlf1 $deXlf2 $deHandbackTime
